$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 0.7
$ws.Range("C2").Value = 61.5
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2
$ws.Range("C3").Value = 61.5
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("C4").Value = 65
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 2
$ws.Range("C5").Value = 59
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 2
$ws.Range("C6").Value = 59
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 3
$ws.Range("C7").Value = 57
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3
$ws.Range("C8").Value = 60
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3
$ws.Range("C9").Value = 62
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3
$ws.Range("C10").Value = 61
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3
$ws.Range("C11").Value = 38
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 7
$ws.Range("C12").Value = 33
$ws.Range("D12").Value = 8
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 8
$ws.Range("A13").Value = 0.8999999999999999
$ws.Range("C13").Value = 67
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 2
$ws.Range("C14").Value = 68
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 2
$ws.Range("C15").Value = 70
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1
$ws.Range("C16").Value = 70
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 2
$ws.Range("C17").Value = 67
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 2
$ws.Range("C18").Value = 65
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 2
$ws.Range("C19").Value = 65
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 4
$ws.Range("C21").Value = 54
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 4
$ws.Range("C22").Value = 41
$ws.Range("D22").Value = 8
$ws.Range("G22").Value = 7
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = 5
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 7
